$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.953.98'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.91%  '
$ws.Range("D3").Value = '''1.867.09'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.18%  '
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").Value = '''318.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.52%  '
$ws.Range("E6").Value = '  -0.28%  '
$ws.Range("D7").Value = '''0.4339'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.94%  '
$ws.Range("D8").Value = '''0.3744'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.34%  '
$ws.Range("D9").Value = '''0.07453'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.20%  '
$ws.Range("D10").Value = '''0.9323'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.85%  '
$ws.Range("D11").Value = '''21.09'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.85%  '
$ws.Range("D12").Value = '''1.905.27'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.14%  '
$ws.Range("D13").Value = '''6.720'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.66%  '
$ws.Range("D14").Value = '''5.421'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.89%  '
$ws.Range("D15").Value = '''0.06860'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.11%  '
$ws.Range("D16").Value = '''1.004'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("D17").Value = '''81.10'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.00%  '
$ws.Range("D18").Value = '''0.000009021'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.51%  '
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").Value = '''15.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.87%  '
$ws.Range("D21").Value = '''27.946.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.12%  '
$ws.Range("D22").Value = '''5.118'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.42%  '
$ws.Range("D23").Value = '''11.01'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("D24").Value = '''2.126.75'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.31%  '
$ws.Range("E25").Value = '  -2.69%  '
$ws.Range("D26").Value = '''152.91'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.33%  '
$ws.Range("D27").Value = '''18.53'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.18%  '
$ws.Range("D28").Value = '''5.529'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.54%  '
$ws.Range("D29").Value = '''113.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.10%  '
$ws.Range("D30").Value = '''1.688'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.51%  '
$ws.Range("D31").Value = '''0.09006'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.65%  '
$ws.Range("D32").Value = '''0.8051'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.01%  '
$ws.Range("D33").Value = '''4.776'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.03%  '
$ws.Range("D34").Value = '''1.176'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.82%  '
$ws.Range("D35").Value = '''2.945'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.72%  '
$ws.Range("D36").Value = '''1.003'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.17%  '
$ws.Range("D37").Value = '''0.05492'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.81%  '
$ws.Range("D38").Value = '''1.117'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.70%  '
$ws.Range("D39").Value = '''0.01973'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.13%  '
$ws.Range("D40").Value = '''2.984'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.40%  '
$ws.Range("D41").Value = '''0.5244'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.98%  '
$ws.Range("D42").Value = '''0.1694'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.71%  '
$ws.Range("D43").Value = '''6.945'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.90%  '
$ws.Range("D44").Value = '''8.745'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.44%  '
$ws.Range("D45").Value = '''0.06739'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.88%  '
$ws.Range("D46").Value = '''0.4873'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.58%  '
$ws.Range("D47").Value = '''10.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.13%  '
$ws.Range("D48").Value = '''106.52'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.57%  '
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("D50").Value = '''1.665'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.02%  '
$ws.Range("D51").Value = '''1.879'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -14.85%  '
